$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value is purely numeric-looking need to be forced to
# Text format first, otherwise Excel auto-converts them to numbers and
# trailing/leading zeros (e.g. '1.330' vs 1.33) get lost.
$textCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D18", "D20", "D21", "D22", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D41", "D43", "D45", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated Price (column D) and Volume(1h) (column E) values
$ws.Range("D2").Value = "26.050.38"
$ws.Range("D3").Value = "1.668.66"
$ws.Range("E3").Value = "  -1.73%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "216.97"
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("D6").Value = "0.5115"
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "0.2658"
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("D9").Value = "0.06415"
$ws.Range("E9").Value = "  +1.86%  "
$ws.Range("D10").Value = "21.88"
$ws.Range("E10").Value = "  -1.57%  "
$ws.Range("D11").Value = "0.07435"
$ws.Range("E11").Value = "  +1.01%  "
$ws.Range("D12").Value = "1.694.39"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").Value = "4.503"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("D14").Value = "0.5837"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").Value = "0.000008579"
$ws.Range("E15").Value = "  +1.36%  "
$ws.Range("D16").Value = "64.39"
$ws.Range("E16").Value = "  -1.87%  "
$ws.Range("D17").Value = "26.096.86"
$ws.Range("E17").Value = "  -2.26%  "
$ws.Range("D18").Value = "4.953"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").Value = "10.78"
$ws.Range("E20").Value = "  -2.09%  "
$ws.Range("D21").Value = "190.64"
$ws.Range("E21").Value = "  +1.81%  "
$ws.Range("D22").Value = "6.238"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "145.07"
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("D25").Value = "7.636"
$ws.Range("E25").Value = "  +1.61%  "
$ws.Range("D26").Value = "0.1202"
$ws.Range("E26").Value = "  +3.37%  "
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("D28").Value = "0.06619"
$ws.Range("E28").Value = "  +16.98%  "
$ws.Range("D29").Value = "1.330"
$ws.Range("E29").Value = "  -1.74%  "
$ws.Range("D30").Value = "1.316"
$ws.Range("E30").Value = "  -1.75%  "
$ws.Range("D31").Value = "3.547"
$ws.Range("E31").Value = "  +1.04%  "
$ws.Range("D32").Value = "3.526"
$ws.Range("E32").Value = "  +1.04%  "
$ws.Range("D33").Value = "1.646"
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("D34").Value = "1.019"
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("D35").Value = "0.6100"
$ws.Range("E35").Value = "  +1.32%  "
$ws.Range("D36").Value = "2.369"
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("D38").Value = "6.259"
$ws.Range("E38").Value = "  +7.01%  "
$ws.Range("D39").Value = "0.01604"
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("D40").Value = "1.086.96"
$ws.Range("E40").Value = "  -1.35%  "
$ws.Range("D41").Value = "0.8619"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("E42").Value = "  +0.60%  "
$ws.Range("D43").Value = "100.39"
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("D44").Value = "1.816.73"
$ws.Range("E44").Value = "  -2.13%  "
$ws.Range("D45").Value = "0.00000000113"
$ws.Range("E45").Value = "  +2.12%  "
$ws.Range("E46").Value = "  -0.75%  "
$ws.Range("D47").Value = "1.011"
$ws.Range("E47").Value = "  +0.64%  "
$ws.Range("D48").Value = "8.042"
$ws.Range("E48").Value = "  -1.72%  "
$ws.Range("D49").Value = "0.05235"
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("D50").Value = "0.4287"
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("D51").Value = "6.004"
$ws.Range("E51").Value = "  +3.74%  "
